$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '57.998.44'
$ws.Range('E2').Value = '  +2.67%  '

$ws.Range('D3').Value = '3.069.95'
$ws.Range('E3').Value = '  +2.53%  '

$ws.Range('E4').Value = '  -0.01%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '520.56'
$ws.Range('E5').Value = '  +2.69%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '142.57'
$ws.Range('E6').Value = '  +3.57%  '

$ws.Range('E7').Value = '  -0.02%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.435'
$ws.Range('E8').Value = '  +1.27%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '7.29'
$ws.Range('E9').Value = '  +1.49%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.108'
$ws.Range('E10').Value = '  +0.06%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.377'
$ws.Range('E11').Value = '  +2.85%  '

$ws.Range('D12').Value = '3.589.03'
$ws.Range('E12').Value = '  +2.73%  '

$ws.Range('E13').Value = '  +3.28%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '25.83'
$ws.Range('E14').Value = '  +0.36%  '

$ws.Range('E15').Value = '  +0.45%  '

$ws.Range('D16').Value = '57.995.89'
$ws.Range('E16').Value = '  +2.82%  '

$ws.Range('D17').Value = '3.065.36'
$ws.Range('E17').Value = '  +2.81%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '6.09'
$ws.Range('E18').Value = '  +1.63%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '12.90'
$ws.Range('E19').Value = '  -0.31%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '8.16'
$ws.Range('E20').Value = '  +1.10%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '330.77'
$ws.Range('E21').Value = '  -0.44%  '

$ws.Range('E22').Value = '  -0.06%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.499'
$ws.Range('E23').Value = '  +1.25%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '65.85'
$ws.Range('E24').Value = '  +1.73%  '

$ws.Range('E25').Value = '  +3.65%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.998'
$ws.Range('E26').Value = '  -0.22%  '

$ws.Range('D27').Value = '0.0₃0902'
$ws.Range('E27').Value = '  -2.10%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '6.41'
$ws.Range('E28').Value = '  +0.61%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '7.22'
$ws.Range('E29').Value = '  +4.26%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.83'
$ws.Range('E30').Value = '  +2.45%  '

$ws.Range('E31').Value = '  +3.17%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '20.69'
$ws.Range('E32').Value = '  +2.09%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '154.87'
$ws.Range('E33').Value = '  +1.29%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.52'
$ws.Range('E34').Value = '  +0.86%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '27.33'
$ws.Range('E35').Value = '  +4.50%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '5.97'
$ws.Range('E36').Value = '  +2.80%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.26'
$ws.Range('E37').Value = '  +0.81%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.0675'
$ws.Range('E38').Value = '  +2.25%  '

$ws.Range('D39').Value = '3.106.76'
$ws.Range('E39').Value = '  +2.84%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '3.93'
$ws.Range('E40').Value = '  +3.37%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '36.73'
$ws.Range('E41').Value = '  -0.55%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.999'

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.653'
$ws.Range('E43').Value = '  +0.35%  '

$ws.Range('D44').Value = '2.275.23'
$ws.Range('E44').Value = '  +4.30%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0259'
$ws.Range('E45').Value = '  +10.25%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '20.91'
$ws.Range('E46').Value = '  +7.05%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.37'
$ws.Range('E47').Value = '  +1.03%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '5.89'
$ws.Range('E48').Value = '  +0.97%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.926'
$ws.Range('E49').Value = '  +0.35%  '

$ws.Range('B50').Value = 'SuiNetwork'
$ws.Range('C50').Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.731'
$ws.Range('E50').Value = '  +8.76%  '

$ws.Range('B51').Value = 'Bittensor'
$ws.Range('C51').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '261.31'
$ws.Range('E51').Value = '  +14.54%  '
